$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): add new headers in M1:O1 ---
$ws.Range("M1").Value = "renewd"
$ws.Range("N1").Value = "PlanID"
$ws.Range("O1").Value = "iteration"

# Copy the header style/format (bold, centered, bordered) from L1 onto the new header cells
$ws.Range("L1").Copy()
$ws.Range("M1:O1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Data rows (2..34): add new data in M:O for each row ---
for ($r = 2; $r -le 34; $r++) {
    $ws.Cells.Item($r, 13).Value = "before"     # Column M
    $ws.Cells.Item($r, 14).Value = 20160319      # Column N
    $ws.Cells.Item($r, 15).Value = 16            # Column O
}
